# Wed, Jun 17, 2020  1:05:42 PM
#
# Switches the deck's applied Design theme from the custom "Integral"
# ("Red Violet" colour scheme) over to the built-in default "Office
# Theme" ("Office" colour scheme). In the PowerPoint UI this is the
# Design tab > Variants > Colors (or simply picking a different Design
# thumbnail) action - it recolours every slide (they all share the one
# Slide Master / theme part) without touching any shape, text or
# layout content.
#
# Helper: build the VBA-style long RGB value (0x00BBGGRR) that
# PowerPoint's ThemeColor.RGB property expects from a normal
# #RRGGBB hex triple.
function ConvertTo-VbaRgb {
    param(
        [int]$r,
        [int]$g,
        [int]$b
    )
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# All slides in this deck hang off the single Slide Master, so grabbing
# the theme color scheme from any one slide gives us the presentation's
# shared theme.
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Target palette: the stock Office theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - the standard 12-slot theme colour order).
$officeColors = @(
    (ConvertTo-VbaRgb 0x00 0x00 0x00),  # 1  dk1      #000000
    (ConvertTo-VbaRgb 0xFF 0xFF 0xFF),  # 2  lt1      #FFFFFF
    (ConvertTo-VbaRgb 0x44 0x54 0x6A),  # 3  dk2      #44546A
    (ConvertTo-VbaRgb 0xE7 0xE6 0xE6),  # 4  lt2      #E7E6E6
    (ConvertTo-VbaRgb 0x5B 0x9B 0xD5),  # 5  accent1  #5B9BD5
    (ConvertTo-VbaRgb 0xED 0x7D 0x31),  # 6  accent2  #ED7D31
    (ConvertTo-VbaRgb 0xA5 0xA5 0xA5),  # 7  accent3  #A5A5A5
    (ConvertTo-VbaRgb 0xFF 0xC0 0x00),  # 8  accent4  #FFC000
    (ConvertTo-VbaRgb 0x44 0x72 0xC4),  # 9  accent5  #4472C4
    (ConvertTo-VbaRgb 0x70 0xAD 0x47),  # 10 accent6  #70AD47
    (ConvertTo-VbaRgb 0x05 0x63 0xC1),  # 11 hlink    #0563C1
    (ConvertTo-VbaRgb 0x95 0x4F 0x72)   # 12 folHlink #954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
